$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-10-02 Wednesday" "2024-10-03 Thursday"

Replace-Text "55×66=3630" "15×61=915"
Replace-Text "95×86=8170" "53×91=4823"
Replace-Text "68×93=6324" "78×44=3432"
Replace-Text "54×51=2754" "69×58=4002"
Replace-Text "30×32=960" "11×67=737"
Replace-Text "21×49=1029" "32×21=672"
Replace-Text "81×66=5346" "25×68=1700"
Replace-Text "72×94=6768" "14×57=798"
Replace-Text "27×38=1026" "11×42=462"
Replace-Text "94×55=5170" "37×53=1961"
Replace-Text "82×11=902" "51×92=4692"
Replace-Text "31×36=1116" "37×84=3108"
Replace-Text "58×74=4292" "96×92=8832"
Replace-Text "56×74=4144" "15×68=1020"
Replace-Text "41×33=1353" "33×50=1650"
Replace-Text "29×83=2407" "88×52=4576"
Replace-Text "64×94=6016" "85×96=8160"
Replace-Text "65×40=2600" "64×13=832"
Replace-Text "44×56=2464" "33×17=561"
Replace-Text "42×68=2856" "46×49=2254"
Replace-Text "66×98=6468" "80×65=5200"
Replace-Text "58×32=1856" "59×25=1475"
Replace-Text "77×73=5621" "68×60=4080"
Replace-Text "76×93=7068" "54×47=2538"
Replace-Text "69×27=1863" "37×74=2738"
